$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5 (Q0-Q3) with new values ---
$ws.Range("B2").Value = -0.07739702482877463
$ws.Range("C2").Value = 1.209708611120026
$ws.Range("D2").Value = 5.87897926970247
$ws.Range("E2").Value = 2.424660650421512
$ws.Range("F2").Value = 2.450204039719076
$ws.Range("G2").Value = 46

$ws.Range("B3").Value = 0.02031028170919459
$ws.Range("C3").Value = 0.9834245770809631
$ws.Range("D3").Value = 5.528573558696921
$ws.Range("E3").Value = 2.351291891428395
$ws.Range("F3").Value = 2.377772292513241
$ws.Range("G3").Value = 45

$ws.Range("B4").Value = 0.0743662420822289
$ws.Range("C4").Value = 1.136542105538119
$ws.Range("D4").Value = 5.807295300627064
$ws.Range("E4").Value = 2.409833044139586
$ws.Range("F4").Value = 2.436532316477507
$ws.Range("G4").Value = 44

$ws.Range("B5").Value = 0.09081030254506482
$ws.Range("C5").Value = 1.050097083832656
$ws.Range("D5").Value = 5.397854852969357
$ws.Range("E5").Value = 2.323328399725135
$ws.Range("F5").Value = 2.349027958552645
$ws.Range("G5").Value = 43

# --- Add new rows 6-11 (Q4-Q9) ---
# First apply the row-label style (bold, bordered, centered) used by A2:A5
# to the new label cells A6:A11, matching the formatting of existing rows.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A6:A11").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A6").Value = "Q4"
$ws.Range("B6").Value = -0.007438971329838501
$ws.Range("C6").Value = 1.151079626825395
$ws.Range("D6").Value = 5.979540328991934
$ws.Range("E6").Value = 2.445309863594374
$ws.Range("F6").Value = 2.474939613007068
$ws.Range("G6").Value = 42

$ws.Range("A7").Value = "Q5"
$ws.Range("B7").Value = 0.1351926068579674
$ws.Range("C7").Value = 1.164807331735655
$ws.Range("D7").Value = 6.217927003399993
$ws.Range("E7").Value = 2.493577150079779
$ws.Range("F7").Value = 2.520841369763717
$ws.Range("G7").Value = 41

$ws.Range("A8").Value = "Q6"
$ws.Range("B8").Value = 0.01744900449909736
$ws.Range("C8").Value = 1.177156600770578
$ws.Range("D8").Value = 6.253386703988633
$ws.Range("E8").Value = 2.500677249064467
$ws.Range("F8").Value = 2.532472641152369
$ws.Range("G8").Value = 40

$ws.Range("A9").Value = "Q7"
$ws.Range("B9").Value = 0.105586766963157
$ws.Range("C9").Value = 1.159725696271973
$ws.Range("D9").Value = 6.386117125364788
$ws.Range("E9").Value = 2.527076794512741
$ws.Range("F9").Value = 2.557876246133412
$ws.Range("G9").Value = 39

$ws.Range("A10").Value = "Q8"
$ws.Range("B10").Value = 0.04157247324465378
$ws.Range("C10").Value = 1.241901496424245
$ws.Range("D10").Value = 6.597015579009581
$ws.Range("E10").Value = 2.568465607908656
$ws.Range("F10").Value = 2.602602220243259
$ws.Range("G10").Value = 38

$ws.Range("A11").Value = "Q9"
$ws.Range("B11").Value = 0.02082192042088303
$ws.Range("C11").Value = 1.103681662169545
$ws.Range("D11").Value = 6.407533224363127
$ws.Range("E11").Value = 2.531310574458047
$ws.Range("F11").Value = 2.566140031814673
$ws.Range("G11").Value = 37
